# Actualización automática 2025-09-08 09:55:08
#
# Applies the refreshed report numbers to the three sheets of the
# "ALMEIDA CUATIN JHONATHANN CARLOS" advisor workbook:
#   1) VENTAS POR GRUPO   - one client's September 240X80 PORCELANATO sale
#   2) VENTA MENSUAL      - mirrors the same sale in the monthly sheet
#   3) CUMPLIMIENTO MENSUAL - full refresh of the per-group compliance table

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Cells.Item(12, 4).Value = 91.58
$wsGrupo.Cells.Item(34, 4).Value = "2 de 32"

# ---------------------------------------------------------------------
# 2) VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Cells.Item(12, 6).Value = 91.58
$wsMensual.Cells.Item(34, 6).Value = 3180

# ---------------------------------------------------------------------
# 3) CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCump = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 - 240X120 PORCELANATO (name unchanged, totals refreshed)
$wsCump.Cells.Item(2, 3).Value = 440.717086537713
$wsCump.Cells.Item(2, 4).Value = 0
$wsCump.Cells.Item(2, 5).Value = 440.717086537713
$wsCump.Cells.Item(2, 6).Value = 0

# Row 3 - 240X80 PORCELANATO (name unchanged, totals refreshed)
$wsCump.Cells.Item(3, 3).Value = 8834.57354940916
$wsCump.Cells.Item(3, 4).Value = 1391.04
$wsCump.Cells.Item(3, 5).Value = 7443.533549409161
$wsCump.Cells.Item(3, 6).Value = 0.1574541195701551

# Row 4 - FREGADEROS DE COCINA (name unchanged, totals refreshed)
$wsCump.Cells.Item(4, 3).Value = 521.61144263264
$wsCump.Cells.Item(4, 4).Value = 145.95
$wsCump.Cells.Item(4, 5).Value = 375.66144263264
$wsCump.Cells.Item(4, 6).Value = 0.2798059783032588

# Row 5 - was GRANITO, now GRIFERIAS
$wsCump.Cells.Item(5, 2).Value = "GRIFERIAS"
$wsCump.Cells.Item(5, 3).Value = 150
$wsCump.Cells.Item(5, 4).Value = 0
$wsCump.Cells.Item(5, 5).Value = 150
$wsCump.Cells.Item(5, 6).Value = 0

# Row 6 - was GRIFERIAS, now INODOROS
$wsCump.Cells.Item(6, 2).Value = "INODOROS"
$wsCump.Cells.Item(6, 3).Value = 814.123430808873
$wsCump.Cells.Item(6, 4).Value = 0
$wsCump.Cells.Item(6, 5).Value = 814.123430808873
$wsCump.Cells.Item(6, 6).Value = 0

# Row 7 - was INODOROS, now LAVABOS
$wsCump.Cells.Item(7, 2).Value = "LAVABOS"
$wsCump.Cells.Item(7, 3).Value = 221.677754071894
$wsCump.Cells.Item(7, 4).Value = 0
$wsCump.Cells.Item(7, 5).Value = 221.677754071894
$wsCump.Cells.Item(7, 6).Value = 0

# Row 8 - was LAVABOS, now NO RESURTIBLES
$wsCump.Cells.Item(8, 2).Value = "NO RESURTIBLES"
$wsCump.Cells.Item(8, 3).Value = 480.217743214072
$wsCump.Cells.Item(8, 4).Value = 0
$wsCump.Cells.Item(8, 5).Value = 480.217743214072
$wsCump.Cells.Item(8, 6).Value = 0

# Row 9 - was LED, now OTROS
$wsCump.Cells.Item(9, 2).Value = "OTROS"
$wsCump.Cells.Item(9, 3).Value = 0
$wsCump.Cells.Item(9, 4).Value = 0
$wsCump.Cells.Item(9, 5).Value = 0
$wsCump.Cells.Item(9, 6).Value = 0

# Row 10 - was NO RESURTIBLES, now PANELES DECORATIVOS
$wsCump.Cells.Item(10, 2).Value = "PANELES DECORATIVOS"
$wsCump.Cells.Item(10, 3).Value = 388.107983534392
$wsCump.Cells.Item(10, 4).Value = 0
$wsCump.Cells.Item(10, 5).Value = 388.107983534392
$wsCump.Cells.Item(10, 6).Value = 0

# Row 11 - was OTROS, now PIEDRA SINTERIZADA
$wsCump.Cells.Item(11, 2).Value = "PIEDRA SINTERIZADA"
$wsCump.Cells.Item(11, 3).Value = 2922.22458185274
$wsCump.Cells.Item(11, 4).Value = 0
$wsCump.Cells.Item(11, 5).Value = 2922.22458185274
$wsCump.Cells.Item(11, 6).Value = 0

# Row 12 - was PANELES DECORATIVOS, now PORCELANATO
$wsCump.Cells.Item(12, 2).Value = "PORCELANATO"
$wsCump.Cells.Item(12, 3).Value = 22433.7553751766
$wsCump.Cells.Item(12, 4).Value = 1643.01
$wsCump.Cells.Item(12, 5).Value = 20790.7453751766
$wsCump.Cells.Item(12, 6).Value = 0.07323829526188128

# Row 13 - was PANELES PU, now PUERTAS DE SEGURIDAD
$wsCump.Cells.Item(13, 2).Value = "PUERTAS DE SEGURIDAD"
$wsCump.Cells.Item(13, 3).Value = 111.043665120341
$wsCump.Cells.Item(13, 4).Value = 0
$wsCump.Cells.Item(13, 5).Value = 111.043665120341
$wsCump.Cells.Item(13, 6).Value = 0

# Row 14 - was PANELES PVC, now SAL SOLUBLE
$wsCump.Cells.Item(14, 2).Value = "SAL SOLUBLE"
$wsCump.Cells.Item(14, 3).Value = 1424.9662010375
$wsCump.Cells.Item(14, 4).Value = 0
$wsCump.Cells.Item(14, 5).Value = 1424.9662010375
$wsCump.Cells.Item(14, 6).Value = 0

# Row 15 used to be PIEDRA SINTERIZADA; it now becomes the TOTAL row
# (the old TOTAL row 19, and the old rows 16-18, are removed below).
$wsCump.Cells.Item(15, 1).ClearContents()
$wsCump.Cells.Item(15, 2).Value = "TOTAL"
$wsCump.Cells.Item(15, 2).HorizontalAlignment = -4152
$wsCump.Cells.Item(15, 3).Value = 38743.01881339593
$wsCump.Cells.Item(15, 4).Value = 3180
$wsCump.Cells.Item(15, 5).Value = 35563.01881339593
$wsCump.Cells.Item(15, 6).Value = 0.08207930350797732

# Remove the now-obsolete rows (old PORCELANATO / PUERTAS DE SEGURIDAD /
# SAL SOLUBLE / TOTAL rows), shrinking the table from 19 to 15 rows.
$wsCump.Rows("16:19").Delete()

# Column widths narrowed slightly for D:F (ColumnWidth getter/setter has a
# fixed +0.83 character offset vs. the stored column width in this engine).
$wsCump.Columns.Item(4).ColumnWidth = 13 - 0.83
$wsCump.Columns.Item(5).ColumnWidth = 23 - 0.83
$wsCump.Columns.Item(6).ColumnWidth = 25 - 0.83
